# Blood Glucose Readings - add "EX2" (second extra/validation) reading columns
# to the Readings sheet: Finger EX2 / Time EX2 / Value EX2, inserted right
# before the existing "Notes" column, plus new sample data for 7/13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Readings")

# --- Insert 3 new columns before column L ("Notes"), pushing Notes to O ---
$ws.Range("L1:N1").EntireColumn.Insert()

# --- New column headers ---
$ws.Range("L1").Value = "Finger EX2"
$ws.Range("M1").Value = "Time EX2"
$ws.Range("N1").Value = "Value EX2"

# The column insert copies column K's formatting into L:N for the rows that
# had a "Value EX" placeholder (rows 2-13). Column L ("Finger EX2") doesn't
# carry a placeholder style/value in the source data, so strip it back out.
$ws.Range("L2:L16").Clear()

# Rows 14-15 need a blank "Value EX2" placeholder (integer format) even
# though column K had nothing there for the insert to copy from; row 16
# should stay completely empty, and no row past 13 should have a "Time EX2"
# placeholder.
$ws.Range("N14:N15").NumberFormat = "0"
$ws.Range("M14:M16").Clear()
$ws.Range("N16").Clear()

# --- New reading data for row 12 (7/13) ---
$ws.Range("G12").Value = 0.79305555555555551
$ws.Range("H12").Value = 69
$ws.Range("L12").Value = 2
$ws.Range("M12").Value = 0.81874999999999998
$ws.Range("M12").NumberFormat = $ws.Range("J12").NumberFormat
$ws.Range("N12").Value = 78

# Notes column (now O) text for row 12 changes to the new plural note
$ws.Range("O12").Value = "Extra readings to validate CGM readings"

# --- Keep the _FilterDatabase defined name in sync with the widened table ---
$fdb = $ws.Names.Item("_xlnm._FilterDatabase")
$fdb.RefersTo = "=Readings!`$A`$1:`$O`$1"
